$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: 질산 납(II) (Lead(II) nitrate) — add missing molecular formula
$ws.Range("E6").Value = "Pb(NO<sub>3</sub>)<sub>2</sub>"

# Row 28: calcium carbonate formula typo fix Ca(CO)3 -> CaCO3
$ws.Range("E28").Value = "CaCO<sub>3</sub>"

# Row 43: Toluidine Blue formula correction
$ws.Range("E43").Value = "2C<sub>15</sub>H<sub>16</sub>N<sub>3</sub>S·Cl·ZnCl<sub>2</sub>"

# Row 46: Magnesium chloride hexahydrate — fix CAS number typo and formula
$ws.Range("A46").Value = '"7791-18-6"'
$ws.Range("E46").Value = "MgCl<sub>2</sub>·6H<sub>2</sub>O"

# Row 54: Sodium bisulfite formula correction
$ws.Range("E54").Value = "NaHSO<sub>3</sub>"
